# Swap the full contents of row 12 and row 13 on the active sheet.
# (The two data-rows were reordered; every cell's value moves from one
# row to the other -- including a few cells that exist on only one side
# of the swap, e.g. K/L/X which row 13 had and row 12 didn't.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = 12
$row2 = 13
$lastCol = 51   # column AY

# Columns whose values look numeric/date/time ("1", "101", "2023-10-02",
# "00:00", ...) but are actually stored as plain text in this sheet. Excel's
# COM layer silently re-types a cell to Number/Date the moment such a string
# is written back into a "General" formatted cell, so these columns must be
# forced to a Text number format before the write to preserve their type.
$textColumns = @(9, 24, 25, 26, 27, 28)   # I, X, Y, Z, AA, AB

for ($c = 1; $c -le $lastCol; $c++) {
    $cell1 = $ws.Cells.Item($row1, $c)
    $cell2 = $ws.Cells.Item($row2, $c)
    $forceText = $textColumns -contains $c

    if ($forceText) {
        # `.Text` always returns the literal display string, never coerced
        # to a number/date/bool the way `.Value2` can be.
        $raw1 = $cell1.Text
        $raw2 = $cell2.Text
        $isBlank1 = ($cell1.Value2 -eq $null)
        $isBlank2 = ($cell2.Value2 -eq $null)
        $old1 = $null; if (-not $isBlank1) { $old1 = $raw1 }
        $old2 = $null; if (-not $isBlank2) { $old2 = $raw2 }
    } else {
        $old1 = $cell1.Value2
        $old2 = $cell2.Value2
    }

    if ($old1 -eq $null -and $old2 -eq $null) {
        continue
    }

    if ($old2 -ne $null) {
        if ($forceText) { $cell1.NumberFormat = "@" }
        $cell1.Value2 = $old2
    } elseif ($old1 -ne $null) {
        $cell1.ClearContents()
    }

    if ($old1 -ne $null) {
        if ($forceText) { $cell2.NumberFormat = "@" }
        $cell2.Value2 = $old1
    } elseif ($old2 -ne $null) {
        $cell2.ClearContents()
    }
}
